$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "0.9970", "30.939.21") keep their exact textual representation
# instead of being parsed into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.939.21'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '1.948.69'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('D4').Value = '0.9971'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '245.18'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').Value = '0.9970'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').Value = '0.4876'
$ws.Range('E7').Value = '  +1.22%  '
$ws.Range('D8').Value = '0.2963'
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('D9').Value = '0.06826'
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').Value = '19.16'
$ws.Range('D11').Value = '107.07'
$ws.Range('E11').Value = '  -4.53%  '
$ws.Range('D12').Value = '1.955.92'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '0.07728'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').Value = '5.474'
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('D15').Value = '0.7060'
$ws.Range('E15').Value = '  +2.79%  '
$ws.Range('D16').Value = '281.71'
$ws.Range('E16').Value = '  -4.22%  '
$ws.Range('D17').Value = '30.953.02'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = '0.000007729'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').Value = '2.206.38'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Value = '0.9962'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').Value = '5.493'
$ws.Range('E22').Value = '  -2.90%  '
$ws.Range('D23').Value = '0.9989'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '6.495'
$ws.Range('E24').Value = '  -1.64%  '
$ws.Range('D25').Value = '9.823'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('D26').Value = '169.03'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '19.94'
$ws.Range('E27').Value = '  -1.68%  '
$ws.Range('D28').Value = '2.214'
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('D29').Value = '0.1052'
$ws.Range('E29').Value = '  -3.24%  '
$ws.Range('E30').Value = '  -1.76%  '
$ws.Range('D31').Value = '1.583'
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').Value = '4.564'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').Value = '4.474'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('D34').Value = '0.04948'
$ws.Range('E34').Value = '  -2.42%  '
$ws.Range('D35').Value = '0.7655'
$ws.Range('E35').Value = '  -1.30%  '
$ws.Range('D36').Value = '1.172'
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('D37').Value = '2.718'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').Value = '0.02024'
$ws.Range('E38').Value = '  -2.47%  '
$ws.Range('D39').Value = '2.691'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '6.544'
$ws.Range('E40').Value = '  +9.16%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '2.157'
$ws.Range('E41').Value = '  +4.58%  '
$ws.Range('D42').Value = '74.95'
$ws.Range('E42').Value = '  +7.77%  '
$ws.Range('D43').Value = '0.4495'
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').Value = '109.37'
$ws.Range('E44').Value = '  -1.41%  '
$ws.Range('D45').Value = '0.8826'
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('D46').Value = '8.167'
$ws.Range('E46').Value = '  +10.47%  '
$ws.Range('D47').Value = '0.9964'
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('D48').Value = '978.91'
$ws.Range('E48').Value = '  +7.35%  '
$ws.Range('D49').Value = '9.408'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('D50').Value = '0.1262'
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('D51').Value = '35.76'
$ws.Range('E51').Value = '  +0.30%  '
